$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 252, shifting existing rows 252:269 down to 253:270.
$ws.Rows("252:252").Insert()

# Populate the newly inserted row 252 with a new weekly price record,
# keeping the same market/category/variety/quality/unit metadata as its
# neighboring rows.
$ws.Range("A252").Value = 11
$ws.Range("B252").Value = "Vega Monumental Concepción"
$ws.Range("C252").Value = "Bíobío"
$ws.Range("D252").Value = 45021
$ws.Range("E252").Value = 8
$ws.Range("F252").Value = 100112003
$ws.Range("G252").Value = "Ajo"
$ws.Range("H252").Value = "Chino"
$ws.Range("I252").Value = "Primera"
$ws.Range("J252").Value = 200
$ws.Range("K252").Value = 14000
$ws.Range("L252").Value = 15000
$ws.Range("M252").Value = 14500
$ws.Range("N252").Value = "`$/caja 10 kilos"
$ws.Range("O252").Value = "China"
$ws.Range("P252").Value = 1450
$ws.Range("Q252").Value = 10
$ws.Range("R252").Value = "Hortaliza"
